# Handback report generation:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The file that was handed off (Latest Handoff File, col C) is now also
#    recorded as the Latest Target File (col E) and Latest Handback File (col F)
#  - Latest Handback DateTime (col G) is stamped with the handback time
# Applies identically to both locale sheets (zh-cn, de-de), each with its own
# handback timestamp.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAt($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            return $h
        }
    }
    return $null
}

function Apply-Handback($ws, $handbackDateTime) {
    # Status -> handed back, in sync with en-US (rows 2 and 3)
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Existing handoff hyperlinks we need to mirror into the target/handback columns
    $hlA2 = Get-HyperlinkAt $ws '$A$2'
    $hlC2 = Get-HyperlinkAt $ws '$C$2'
    $hlA3 = Get-HyperlinkAt $ws '$A$3'
    $hlC3 = Get-HyperlinkAt $ws '$C$3'

    # Row 2: Latest Target File (E) / Latest Handback File (F) mirror the
    # handoff file (A/C), since the handback is in sync with the handoff.
    $ws.Hyperlinks.Add($ws.Range("E2"), $hlA2.Address(), [Type]::Missing, [Type]::Missing, $hlA2.TextToDisplay())
    $ws.Hyperlinks.Add($ws.Range("F2"), $hlC2.Address(), [Type]::Missing, [Type]::Missing, $hlC2.TextToDisplay())

    # Row 3
    $ws.Hyperlinks.Add($ws.Range("E3"), $hlA3.Address(), [Type]::Missing, [Type]::Missing, $hlA3.TextToDisplay())
    $ws.Hyperlinks.Add($ws.Range("F3"), $hlC3.Address(), [Type]::Missing, [Type]::Missing, $hlC3.TextToDisplay())

    # Latest Handback DateTime (G) for rows 2 and 3
    $ws.Range("G2").Value = $handbackDateTime
    $ws.Range("G3").Value = $handbackDateTime
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Apply-Handback $wsZhCn "2016-02-17 06:43:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
Apply-Handback $wsDeDe "2016-02-17 06:43:34"

# The Overview sheet's status cells share the same underlying text as the
# per-locale "Status" column (col B/C on Overview mirror "Ready for handoff").
# Keep them in sync with the new status text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
